# Add the new "Typ" column (P) with crs/grp type markers, renumber the
# leading Kurs-/Gruppen-Ref-Ids (col A) and the Hierarchie-Id (col D), and
# drop the now-obsolete Kurs-Ref-Id in row 4 (new/“create” row has none).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New "Typ" header (row 1, bold/red "muss ausgefüllt"-style header cell,
#     with the same date-ish number format column P already carried) ---
$ws.Range("P1").Value = "Typ"
$ws.Range("P1").Font.Color = 255
$ws.Range("P1").Interior.ThemeColor = 2
$ws.Range("P1").Interior.TintAndShade = 0

# --- Row 2 (German sub-header) gets the same label ---
$ws.Range("P2").Value = "Typ"

# --- Data rows: Typ column values ---
$ws.Range("P3").Value = "crs"
$ws.Range("P4").Value = "grp"
$ws.Range("P5").Value = "crs"
$ws.Range("P6").Value = "crs"
$ws.Range("P7").Value = "crs"
$ws.Range("P8").Value = "crs"

# --- Renumbered reference ids in column A ---
$ws.Range("A3").Value = 83
$ws.Range("A5").Value = 85
$ws.Range("A6").Value = 86
$ws.Range("A7").Value = 87
$ws.Range("A8").Value = 88

# Row 4 no longer carries a Kurs-Ref-Id (this row demonstrates creating a
# brand-new group, so the ref-id cell is cleared entirely).
$ws.Range("A4").ClearContents()

# --- Renumbered hierarchy ids in column D ---
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 83
$ws.Range("D5").Value = 83
$ws.Range("D6").Value = 83
$ws.Range("D7").Value = 83
$ws.Range("D8").Value = 83

# --- Move the active selection cursor the way the author left it ---
$ws.Range("C4").Select()
